$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-20 Sunday" "2025-07-21 Monday"

Replace-Text "49×71=" "31×61="
Replace-Text "71×70=" "27×92="
Replace-Text "36×55=" "74×64="
Replace-Text "19×63=" "50×24="
Replace-Text "23×79=" "43×92="

Replace-Text "65×17=" "70×24="
Replace-Text "36×86=" "27×21="
Replace-Text "68×85=" "43×78="
Replace-Text "83×89=" "11×82="
Replace-Text "21×70=" "97×23="

Replace-Text "65×13=" "47×98="
Replace-Text "46×29=" "23×80="
Replace-Text "80×30=" "92×61="
Replace-Text "17×57=" "63×99="
Replace-Text "21×61=" "83×63="

Replace-Text "19×20=" "72×80="
Replace-Text "23×32=" "62×64="
Replace-Text "34×67=" "17×14="
Replace-Text "33×30=" "56×30="
Replace-Text "92×57=" "26×14="

Replace-Text "97×32=" "50×44="
Replace-Text "95×55=" "11×47="
Replace-Text "41×58=" "65×13="
Replace-Text "83×24=" "75×66="
Replace-Text "15×53=" "84×94="
